# Debug: Fixing the errors
# Add the new "Soil_Dynamic_Temperature_Factor" parameter row to the
# Default_Inputs worksheet (new row 35, mirrors the existing Parameter/value
# rows above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Soil_Dynamic_Temperature_Factor"
$ws.Range("B35").Value = 4

# Move the selection/active cell onto the newly added row, matching the
# author's saved view state.
[void]$ws.Range("A35").Select()
